$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3-Year Summary")

# Rows 3 & 4: the category labels were swapped (Hardware <-> Software Licenses).
# The SUMIF/sum formulas already reference their own row via relative refs
# (A3/A4), so they keep working correctly once the label text changes.
$ws.Range("A3").Value = "Software Licenses"
$ws.Range("A4").Value = "Hardware"

# Row 5 (Support & Maintenance) is unchanged.

# Row 6 used to be the "TOTAL" (SUM) row; it becomes the "Operational
# Efficiency" category row, pulling its numbers via SUMIF from the
# Infrastructure Costs / Credits sheets, same pattern as rows 3-5 & 7.
$ws.Range("A6").Value = "Operational Efficiency"
$ws.Range("B6").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A6,'Infrastructure Costs'!`$G:`$G)"
$ws.Range("C6").Formula = "=SUMIF(Credits!`$A:`$A,A6,Credits!`$C:`$C)"
$ws.Range("D6").Formula = "=B6+C6"
$ws.Range("E6").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A6,'Infrastructure Costs'!`$H:`$H)"
$ws.Range("F6").Formula = "=SUMIF('Infrastructure Costs'!`$A:`$A,A6,'Infrastructure Costs'!`$I:`$I)"
$ws.Range("G6").Formula = "=D6+E6+F6"

# Row 7: relabel "Cloud Cost Avoidance (Annual)" -> "Cloud Cost Avoidance"
# (matches the Credits sheet's category name exactly); formulas unchanged.
$ws.Range("A7").Value = "Cloud Cost Avoidance"

# Row 8 used to be "Net Investment After Savings" (SUMIF row); it becomes the
# "TOTAL" row, now summing rows 3-7 (categories + Operational Efficiency +
# Cloud Cost Avoidance) instead of doing its own SUMIF lookups.
$ws.Range("A8").Value = "TOTAL"
$ws.Range("B8").Formula = "=SUM(B3:B7)"
$ws.Range("C8").Formula = "=SUM(C3:C7)"
$ws.Range("D8").Formula = "=SUM(D3:D7)"
$ws.Range("E8").Formula = "=SUM(E3:E7)"
$ws.Range("F8").Formula = "=SUM(F3:F7)"
$ws.Range("G8").Formula = "=SUM(G3:G7)"

# A blank row 9 now exists below the (moved) TOTAL row; touch then clear it
# so the sheet's used range grows to A1:G9, matching the stray-row cleanup
# described in the commit message, while leaving it with no visible content
# or formatting.
$ws.Range("A9").NumberFormat = "General"
$ws.Range("A9").ClearFormats()
